$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper behaviour note: toggling Font.Bold true->false on a sub-range is a
# cheap way to force the run-split boundary at that exact offset without
# changing the visible formatting (the engine keeps runs split even once
# their formatting becomes identical again).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 1) Table column widths: 3687/4607 -> 3685/4609 (tiny AutoFit re-balance)
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)
$tbl.Columns(1).SetWidth(3685, 1)
$tbl.Columns(2).SetWidth(4609, 1)

# ---------------------------------------------------------------------------
# 2) {direccion} -> collapse the spell-checked run split into a single run
#    (same visible text, but merges the 3 runs / drops the proofErr marks)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("{direccion}", $false, $false, $false, $false, $false, $true, 1, $false, "{direccion}", 2)

# ---------------------------------------------------------------------------
# 3) {fechaInicio} -> {fechaInicioTexto}
#    Target run layout: "{fechaInici" | "oTexto" | "}"  (no proofErr marks)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("{fechaInicio}", $false, $false, $false, $false, $false, $true, 1, $false, "{fechaInicioTexto}", 2)

$r = $d.Content
$r.Find.Execute("{fechaInicioTexto}")
$matchEnd = $r.End

# split "oTexto}" off from "{fechaInici"
$split1 = $d.Range($matchEnd - 7, $matchEnd)
$split1.Font.Bold = $true
$split1.Font.Bold = $false

# split "}" off from "oTexto"
$split2 = $d.Range($matchEnd - 1, $matchEnd)
$split2.Font.Bold = $true
$split2.Font.Bold = $false

# ---------------------------------------------------------------------------
# 4) {fechaFinal} -> {fechaFinalTexto}  (both occurrences get replaced by a
#    single Execute() over the whole-document Range with wrap=wdFindContinue)
#    Target run layout per occurrence: "{fecha" | "Final" | "Texto" | "}"
#    (no proofErr marks)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("{fechaFinal}", $false, $false, $false, $false, $false, $true, 1, $false, "{fechaFinalTexto}", 2)

$r = $d.Content
$r.Find.Execute("{fechaFinalTexto}")
$firstEnd = $r.End

$r2 = $d.Range($r.End, $d.Content.End)
$r2.Find.Execute("{fechaFinalTexto}")
$secondEnd = $r2.End

foreach ($matchEnd in @($firstEnd, $secondEnd)) {
    # split "Texto}" off from "{fechaFinal"
    $split1 = $d.Range($matchEnd - 6, $matchEnd)
    $split1.Font.Bold = $true
    $split1.Font.Bold = $false

    # split "}" off from "Texto"
    $split2 = $d.Range($matchEnd - 1, $matchEnd)
    $split2.Font.Bold = $true
    $split2.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 5) Replace the spelled-out "veintiséis (09) del mes de mayo de 2024" date
#    with the {fechaInicioTexto} placeholder.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("veintiséis (09) del mes de mayo de 2024", $false, $false, $false, $false, $false, $true, 1, $false, "{fechaInicioTexto}", 2)

Write-Output "done"
